$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Formula = "'2014"
$ws.Range("C2").Value2 = 65830.0
$ws.Range("C3").Value2 = 0.0000002194
$ws.Range("C4").Value2 = 25690.0
$ws.Range("C5").Value2 = 20710.0
$ws.Range("C6").Value2 = 4980.0
$ws.Range("C7").Value2 = 3520.0
$ws.Range("C8").Value2 = 1460.0
$ws.Range("C9").Value2 = 0.0000002476
$ws.Range("C10").Value2 = 40140.0
$ws.Range("C11").Value2 = 0.0000002018
$ws.Range("C13").Value2 = 23810.0
$ws.Range("C14").Value2 = 9830.0
$ws.Range("C15").Value2 = 13980.0
$ws.Range("C16").Value2 = 0.0000002607
$ws.Range("C18").Value2 = 281.0
$ws.Range("C19").Value2 = -281.0
$ws.Range("C20").Value2 = 570.0
$ws.Range("C21").Value2 = 746.0
$ws.Range("C22").ClearContents() | Out-Null
$ws.Range("C23").Value2 = 101.0
$ws.Range("C24").Value2 = 0.0000001923
$ws.Range("C25").Value2 = 101.0
$ws.Range("C26").ClearContents() | Out-Null
$ws.Range("C27").Value2 = 17260.0
$ws.Range("C28").Value2 = 0.0000002289
$ws.Range("C30").Value2 = 3640.0
$ws.Range("C31").Value2 = 2870.0
$ws.Range("C32").Value2 = 774.0
$ws.Range("C33").Value2 = 35.0
$ws.Range("C34").Value2 = -43.0
$ws.Range("C38").Value2 = 13620.0
$ws.Range("C40").Value2 = 13620.0
$ws.Range("C41").Value2 = -0.0000003499
$ws.Range("C43").Value2 = 516.0
$ws.Range("C44").Value2 = 740.0
$ws.Range("C46").Value2 = -224.0
$ws.Range("C47").Value2 = 14140.0
$ws.Range("C49").Value2 = 13400.0
$ws.Range("C50").Value2 = 0.00001982
$ws.Range("C51").Value2 = -0.0000003547
$ws.Range("C52").Value2 = 675.94
$ws.Range("C53").Value2 = 0.0000195
$ws.Range("C54").Value2 = -0.0000003544
$ws.Range("C55").Value2 = 687.07
$ws.Range("C56").Value2 = 21300.0
$ws.Range("C57").Value2 = 0.0000002735
$ws.Range("C59").Value2 = 64400.00000000001
$ws.Range("C60").Value2 = 18350.0
$ws.Range("C61").Value2 = 46050.0
$ws.Range("C62").Value2 = 0.0000001347
$ws.Range("C63").Value2 = 0.0000004985
$ws.Range("C64").Value2 = 10850.0
$ws.Range("C65").Value2 = 9380.0
$ws.Range("C66").Value2 = 9610.0
$ws.Range("C67").Value2 = -225.0
$ws.Range("C68").Value2 = 1470.0
$ws.Range("C69").Value2 = 0.0000002821
$ws.Range("C70").Value2 = 0.00000607
$ws.Range("C71").ClearContents() | Out-Null
$ws.Range("C72").ClearContents() | Out-Null
$ws.Range("C76").Value2 = 3410.0
$ws.Range("C77").Value2 = 3410.0
$ws.Range("C78").Value2 = 78660.0
$ws.Range("C79").Value2 = 23880.0
$ws.Range("C80").Value2 = 32750.0
$ws.Range("C81").Value2 = 13330.0
$ws.Range("C83").Value2 = 10920.0
$ws.Range("C84").Value2 = 1950.0
$ws.Range("C85").Value2 = 8860.0
$ws.Range("C86").Value2 = 3080.0
$ws.Range("C87").ClearContents() | Out-Null
$ws.Range("C88").Value2 = 1330.0
$ws.Range("C89").Value2 = 20210.0
$ws.Range("C90").Value2 = 15600.0
$ws.Range("C91").Value2 = 4610.0
$ws.Range("C92").Value2 = 1860.0
$ws.Range("C93").Value2 = 1860.0
$ws.Range("C94").Value2 = 129190.0
$ws.Range("C95").Value2 = 0.0000001415
$ws.Range("C96").Value2 = 2010.0
$ws.Range("C97").Value2 = 2000.0
$ws.Range("C98").Value2 = 10.0
$ws.Range("C99").Value2 = 1720.0
$ws.Range("C100").Value2 = 0.000000537
$ws.Range("C101").Value2 = 96.0
$ws.Range("C102").Value2 = 12960.0
$ws.Range("C104").Value2 = 3070.0
$ws.Range("C105").Value2 = 9890.0
$ws.Range("C106").Value2 = 16780.0
$ws.Range("C107").Value2 = 3230.0
$ws.Range("C108").Value2 = 2990.0
$ws.Range("C109").Value2 = 2990.0
$ws.Range("C111").Value2 = 236.0
$ws.Range("C113").Value2 = 582.0
$ws.Range("C114").Value2 = 758.0
$ws.Range("C115").Value2 = 176.0
$ws.Range("C116").Value2 = 4560.0
$ws.Range("C117").Value2 = 4460.0
$ws.Range("C118").Value2 = 104.0
$ws.Range("C119").Value2 = 25330.0
$ws.Range("C121").Value2 = 0.000000196
$ws.Range("C125").Value2 = 103860.0
$ws.Range("C126").Value2 = 0.68
$ws.Range("C127").Value2 = 75070.0
$ws.Range("C129").Value2 = -980.0
$ws.Range("C130").Value2 = 421.0
$ws.Range("C133").Value2 = 0.000000804
$ws.Range("C134").Value2 = 103860.0
$ws.Range("C135").Value2 = 0.000000804
$ws.Range("C137").Value2 = 103860.0
$ws.Range("C138").Value2 = 129190.0
$ws.Range("C139").Value2 = 14140.0
$ws.Range("C140").Value2 = -0.0000003499
$ws.Range("C141").Value2 = 4980.0
$ws.Range("C142").Value2 = 3520.0
$ws.Range("C143").Value2 = 1460.0
$ws.Range("C144").Value2 = -104.0
$ws.Range("C145").Value2 = -104.0
$ws.Range("C147").Value2 = 2690.0
$ws.Range("C148").Value2 = 21700.0
$ws.Range("C150").Value2 = 672.0
$ws.Range("C151").Value2 = -1640.0
$ws.Range("C152").Value2 = 436.0
$ws.Range("C153").Value2 = 284.0
$ws.Range("C154").Value2 = 22380.0
$ws.Range("C155").Value2 = 0.0000003562
$ws.Range("C156").Value2 = 0.0000003399
$ws.Range("C157").Value2 = -10960.0
$ws.Range("C158").Value2 = -10960.0
$ws.Range("C160").Value2 = 0.00000009210000000000001
$ws.Range("C161").Value2 = -0.0000001665
$ws.Range("C162").Value2 = -4890.0
$ws.Range("C163").Value2 = 386.0
$ws.Range("C164").Value2 = -7000.0
$ws.Range("C165").Value2 = -58310.0
$ws.Range("C166").Value2 = 51320.0
$ws.Range("C168").Value2 = 1400.0
$ws.Range("C169").Value2 = -21060.0
$ws.Range("C170").Value2 = -0.0000001261
$ws.Range("C171").Value2 = -0.0000003198
$ws.Range("C175").ClearContents() | Out-Null
$ws.Range("C176").ClearContents() | Out-Null
$ws.Range("C180").Value2 = -18.0
$ws.Range("C182").Value2 = -18.0
$ws.Range("C183").Value2 = 11630.0
$ws.Range("C184").Value2 = -11640.0
$ws.Range("C185").Value2 = -1420.0
$ws.Range("C186").Value2 = -2070.0
$ws.Range("C187").Value2 = 648.0
$ws.Range("C188").Value2 = -1440.0
$ws.Range("C189").Value2 = -0.0000019361
$ws.Range("C190").Value2 = -0.0000000219
$ws.Range("C191").Value2 = -433.0
$ws.Range("C193").Value2 = -551.0
$ws.Range("C194").Value2 = 11420.0
$ws.Range("C195").Value2 = 0.0000004559000000000001
